$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 245, shifting existing data (rows 245-372) down to (247-374)
$ws.Rows("245:246").Insert()

# Row 245 - new data
$ws.Cells.Item(245, 1).Value = 6
$ws.Cells.Item(245, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(245, 3).Value = "Metropolitana"
$ws.Cells.Item(245, 4).Value = 44455
$ws.Cells.Item(245, 5).Value = 13
$ws.Cells.Item(245, 6).Value = 100112017
$ws.Cells.Item(245, 7).Value = "Apio"
$ws.Cells.Item(245, 8).Value = "Americana (o)"
$ws.Cells.Item(245, 9).Value = "Primera"
$ws.Cells.Item(245, 10).Value = 1700
$ws.Cells.Item(245, 11).Value = 6000
$ws.Cells.Item(245, 12).Value = 7000
$ws.Cells.Item(245, 13).Value = 6471
$ws.Cells.Item(245, 14).Value = "`$/docena de matas"
$ws.Cells.Item(245, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(245, 16).Value = 1078
$ws.Cells.Item(245, 17).Value = 6
$ws.Cells.Item(245, 18).Value = "Hortaliza"

# Row 246 - new data
$ws.Cells.Item(246, 1).Value = 6
$ws.Cells.Item(246, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(246, 3).Value = "Metropolitana"
$ws.Cells.Item(246, 4).Value = 44455
$ws.Cells.Item(246, 5).Value = 13
$ws.Cells.Item(246, 6).Value = 100112017
$ws.Cells.Item(246, 7).Value = "Apio"
$ws.Cells.Item(246, 8).Value = "Americana (o)"
$ws.Cells.Item(246, 9).Value = "Segunda"
$ws.Cells.Item(246, 10).Value = 600
$ws.Cells.Item(246, 11).Value = 5000
$ws.Cells.Item(246, 12).Value = 5000
$ws.Cells.Item(246, 13).Value = 5000
$ws.Cells.Item(246, 14).Value = "`$/docena de matas"
$ws.Cells.Item(246, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(246, 16).Value = 833
$ws.Cells.Item(246, 17).Value = 6
$ws.Cells.Item(246, 18).Value = "Hortaliza"

# Ensure the date style (numFmt) carries correctly for new D cells
$ws.Cells.Item(245, 4).NumberFormat = $ws.Cells.Item(247, 4).NumberFormat
$ws.Cells.Item(246, 4).NumberFormat = $ws.Cells.Item(247, 4).NumberFormat

Write-Host "Done. UsedRange rows: $($ws.UsedRange.Rows.Count)"
